# Disaggregation of commodity Copper
#
# Every year-sheet (2000..2100) lists three generation-technology rows in
# C5:C7 / E5:E7 (and D/F/G, which are always 0):
#   row 5 -> "Photovoltaic plants"
#   row 6 -> "Onshore wind plants"
#   row 7 -> "Offshore wind plants"
#
# The edit swaps the "Photovoltaic plants" and "Onshore wind plants" rows
# (both their label and their numeric value) on every sheet, while the
# "Offshore wind plants" row (7) is left untouched.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $label5 = $ws.Range("C5").Value2
    $label6 = $ws.Range("C6").Value2
    $value5 = $ws.Range("E5").Value2
    $value6 = $ws.Range("E6").Value2

    $ws.Range("C5").Value = $label6
    $ws.Range("C6").Value = $label5
    $ws.Range("E5").Value = $value6
    $ws.Range("E6").Value = $value5
}
